$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new booking rows (7 and 8) for SNOW-744415, matching the existing
# "Family Ski Package" style rows (5/6) where text-like columns (dates,
# phone numbers, empty Special Requests) must stay stored as text rather
# than being coerced to numbers/dates.

$rows = @(
    @{ Row = 7; A = "SNOW-744415"; B = "2026-03-19"; C = "march"; D = "m@b.com"; E = "1213"; F = 2; G = "Ski Adventure"; H = 18500; I = 18500; J = "Confirmed"; K = "2026-02-17" },
    @{ Row = 8; A = "SNOW-744415"; B = "2026-03-20"; C = "march"; D = "m@b.com"; E = "1213"; F = 2; G = "Ski Adventure"; H = 18500; I = 18500; J = "Confirmed"; K = "2026-02-17" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A

    # Text columns that look like dates/numbers: force text format first
    # so Excel doesn't auto-convert them to a date serial / number.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r.B

    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D

    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $r.E

    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J

    $ws.Range("K$row").NumberFormat = "@"
    $ws.Range("K$row").Value = $r.K

    # Special Requests: empty text value (leading apostrophe forces a
    # text-typed empty cell instead of clearing it to blank/number).
    $ws.Range("L$row").Value = "'"
}
